# chore(results): update lottery results 2025-09-20T17:39:43Z
#
# Append the new Pick 4 results row (row 4) to the "Results" sheet.
# All values in this sheet are stored as text (dates like "2025-09-20" and
# phase codes like "250920" must NOT be auto-converted to a date serial or a
# number by Excel), so each cell's number format is forced to Text ("@")
# before the value is written, then the formatting is cleared again so the
# new row keeps the sheet's default (unstyled) look - matching the style of
# the pre-existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A4:E4")
$newRow.NumberFormat = "@"

$ws.Range("A4").Value = "2025-09-20"
$ws.Range("B4").Value = "Pick 4"
$ws.Range("C4").Value = "250920"
$ws.Range("D4").Value = "9-5-3-7"
$ws.Range("E4").Value = "2025-09-20T21:39:43.152+04:00"

# Drop the explicit Text number-format again so row 4 ends up with the same
# (default) style as the existing rows - only the stored value type (text)
# needs to stick, not a lingering "@" format.
$newRow.ClearFormats()
